$d = $word.ActiveDocument

# Change 1: merge the two runs "Розв'язок завдання №" + "2" that make up the
# "Розв'язок завдання №2" heading (bookmark _Toc132872427) into a single run.
# Scope the Find/Replace to that bookmark's range so only this heading
# (and not the identical-looking TOC hyperlink entry) is touched.
$bm = $d.Bookmarks("_Toc132872427")
$r = $bm.Range
$r.Find.Execute("Розв’язок завдання №2", $false, $false, $false, $false, `
                $false, $true, 1, $false, "Розв’язок завдання №2", 2) | Out-Null

# Change 2: the "Результати роботи завдання №2" paragraph switches from
# direct formatting (spacing/indent) to the "Заголовок-2" (pStyle "-20")
# paragraph style.
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -eq "Результати роботи завдання №2`r") {
        $p.Style = "Заголовок-2"
    }
}
